# Grid Bat Cap Growth per Unit Net Revenue.xlsx
# "Draft grid battery edits and fix to battery cost calculation"
#
# Semantic changes captured by this script:
#   1. Fix the battery cost calculation input: GBCGpUNR!B2 (MW per $/MWh net
#      revenue, "All years") goes from 2000 to 750.
#   2. Leave the workbook in the UI state the author last saved in: cursor
#      moved to B3 on the GBCGpUNR sheet (next to the edited value), with the
#      About sheet left as the active/selected tab.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("GBCGpUNR")

# --- fix to battery cost calculation -------------------------------------
$wsData.Range("B2").Value = 750

# --- draft grid battery edits: final selection / active-tab state --------
$wsData.Activate() | Out-Null
$wsData.Range("B3").Select() | Out-Null

$wsAbout.Activate() | Out-Null
